$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value = 44267.55206759619
$ws.Range("D16:D29").Value = 44267.53069166667
$ws.Range("D30:D43").Value = 44267.50927886574
